$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 7952979.5
$ws.Range("I6").Value = 13891464
$ws.Range("J6").Value = 35000
$ws.Range("K6").Value = 41674392
$ws.Range("L6").Value = 105000
$ws.Range("M6").Value = -41674280
$ws.Range("N6").Value = -105224
$ws.Range("H70").Value = 1591.6666
$ws.Range("I70").Value = 1050
$ws.Range("J70").Value = 1788.6364
$ws.Range("K70").Value = 3150
$ws.Range("L70").Value = 5365.9092
$ws.Range("M70").Value = -2880
$ws.Range("N70").Value = -5905.9092
$ws.Range("H73").Value = 1591.6666
$ws.Range("I73").Value = 1050
$ws.Range("J73").Value = 1788.6364
$ws.Range("K73").Value = 3150
$ws.Range("L73").Value = 5365.9092
$ws.Range("M73").Value = -2214
$ws.Range("N73").Value = -7237.9092
$ws.Range("H98").Value = 32039.44
$ws.Range("I98").Value = 652.8095
$ws.Range("J98").Value = 196819.25
$ws.Range("K98").Value = 652.8095
$ws.Range("L98").Value = 196819.25
$ws.Range("M98").Value = 845.1905
$ws.Range("N98").Value = -199815.25
$ws.Range("H122").Value = 32039.44
$ws.Range("I122").Value = 652.8095
$ws.Range("J122").Value = 196819.25
$ws.Range("K122").Value = 1958.4285
$ws.Range("L122").Value = 590457.75
$ws.Range("M122").Value = 491.5715
$ws.Range("N122").Value = -595357.75
$ws.Range("H129").Value = 1059.5714
$ws.Range("I129").Value = 1007.2222
$ws.Range("J129").Value = 1089.9678
$ws.Range("K129").Value = 3021.6666
$ws.Range("L129").Value = 3269.9034
$ws.Range("M129").Value = 1978.3334
$ws.Range("N129").Value = -13269.9034
$ws.Range("H137").Value = 2657275
$ws.Range("I137").Value = 5918082.5
$ws.Range("J137").Value = 7868.875
$ws.Range("K137").Value = 17754247.5
$ws.Range("L137").Value = 23606.625
$ws.Range("M137").Value = -17751697.5
$ws.Range("N137").Value = -28706.625
$ws.Range("H138").Value = 1443.13
$ws.Range("I138").Value = 734.84375
$ws.Range("J138").Value = 1776.4412
$ws.Range("K138").Value = 2204.53125
$ws.Range("L138").Value = 5329.3236
$ws.Range("M138").Value = 2935.46875
$ws.Range("N138").Value = -15609.3236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2063.0557
$ws.Range("I110").Value = 2137.2856
$ws.Range("J110").Value = 1803.25
$ws.Range("K110").Value = 2137.2856
$ws.Range("L110").Value = 1803.25
$ws.Range("M110").Value = -92.28560000000016
$ws.Range("N110").Value = -5893.25
$ws.Range("H132").Value = 11113288
$ws.Range("I132").Value = 16130510
$ws.Range("J132").Value = 3726.0715
$ws.Range("K132").Value = 48391530
$ws.Range("L132").Value = 11178.2145
$ws.Range("M132").Value = -48389000
$ws.Range("N132").Value = -16238.2145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2143.2559
$ws.Range("J99").Value = 2290
$ws.Range("L99").Value = 2290
$ws.Range("N99").Value = -5286
$ws.Range("H105").Value = 3857.0625
$ws.Range("I105").Value = 3145.7144
$ws.Range("J105").Value = 4410.3335
$ws.Range("K105").Value = 3145.7144
$ws.Range("L105").Value = 4410.3335
$ws.Range("M105").Value = -1398.7144
$ws.Range("N105").Value = -7904.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2005.5
$ws.Range("I16").Value = 3011
$ws.Range("K16").Value = 3011
$ws.Range("M16").Value = -2724
$ws.Range("H105").Value = 3717.077
$ws.Range("I105").Value = 3950
$ws.Range("K105").Value = 3950
$ws.Range("M105").Value = -2203
$ws.Range("H113").Value = 2005.5
$ws.Range("I113").Value = 3011
$ws.Range("K113").Value = 3011
$ws.Range("M113").Value = -841
$ws.Range("H124").Value = 39999.668
$ws.Range("J124").Value = 39999.668
$ws.Range("L124").Value = 39999.668
$ws.Range("N124").Value = -44909.668
$ws.Range("H132").Value = 111139.16
$ws.Range("I132").Value = 2602
$ws.Range("K132").Value = 7806
$ws.Range("M132").Value = -5276
$ws.Range("H134").Value = 913644.25
$ws.Range("I134").Value = 1175745.1
$ws.Range("K134").Value = 3527235.3
$ws.Range("M134").Value = -3524700.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1877.5
$ws.Range("I7").Value = 2617.5
$ws.Range("K7").Value = 7852.5
$ws.Range("M7").Value = -7740.5
$ws.Range("H41").Value = 626.25
$ws.Range("I41").Value = 301
$ws.Range("K41").Value = 903
$ws.Range("M41").Value = -565
$ws.Range("H80").Value = 46394896
$ws.Range("J80").Value = 52926880
$ws.Range("L80").Value = 158780640
$ws.Range("N80").Value = -158782512
$ws.Range("H83").Value = 46394896
$ws.Range("J83").Value = 52926880
$ws.Range("L83").Value = 476341920
$ws.Range("N83").Value = -476351280
$ws.Range("H92").Value = 1301.1428
$ws.Range("I92").Value = 1020.4
$ws.Range("J92").Value = 2003
$ws.Range("K92").Value = 3061.2
$ws.Range("L92").Value = 6009
$ws.Range("M92").Value = -1813.2
$ws.Range("N92").Value = -8505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 27131
$ws.Range("J75").Value = 27131
$ws.Range("L75").Value = 27131
$ws.Range("N75").Value = -28879
$ws.Range("H78").Value = 27131
$ws.Range("J78").Value = 27131
$ws.Range("L78").Value = 81393
$ws.Range("N78").Value = -90129
$ws.Range("H94").Value = 25714.285
$ws.Range("J94").Value = 25714.285
$ws.Range("L94").Value = 25714.285
$ws.Range("N94").Value = -27066.285
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
$ws.Range("H141").Value = 45174.668
$ws.Range("J141").Value = 45174.668
$ws.Range("L141").Value = 45174.668
$ws.Range("N141").Value = -55534.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 549.069
$ws.Range("I55").Value = 505.13635
$ws.Range("J55").Value = 687.1429000000001
$ws.Range("K55").Value = 505.13635
$ws.Range("L55").Value = 687.1429000000001
$ws.Range("M55").Value = -332.13635
$ws.Range("N55").Value = -1033.1429
$ws.Range("H68").Value = 3157.7058
$ws.Range("I68").Value = 2898.4167
$ws.Range("K68").Value = 2898.4167
$ws.Range("M68").Value = -2149.4167
$ws.Range("H71").Value = 3157.7058
$ws.Range("I71").Value = 2898.4167
$ws.Range("K71").Value = 14492.0835
$ws.Range("M71").Value = -10748.0835
$ws.Range("H74").Value = 40108.5
$ws.Range("I74").Value = 40000
$ws.Range("J74").Value = 40217
$ws.Range("K74").Value = 40000
$ws.Range("L74").Value = 40217
$ws.Range("M74").Value = -39002
$ws.Range("N74").Value = -42213
$ws.Range("H77").Value = 40108.5
$ws.Range("I77").Value = 40000
$ws.Range("J77").Value = 40217
$ws.Range("K77").Value = 120000
$ws.Range("L77").Value = 120651
$ws.Range("M77").Value = -115008
$ws.Range("N77").Value = -130635
$ws.Range("H82").Value = 6411063
$ws.Range("I82").Value = 994
$ws.Range("J82").Value = 8334084
$ws.Range("K82").Value = 994
$ws.Range("L82").Value = 8334084
$ws.Range("M82").Value = -633
$ws.Range("N82").Value = -8334806
$ws.Range("H85").Value = 6411063
$ws.Range("I85").Value = 994
$ws.Range("J85").Value = 8334084
$ws.Range("K85").Value = 994
$ws.Range("L85").Value = 8334084
$ws.Range("M85").Value = 254
$ws.Range("N85").Value = -8336580
$ws.Range("H132").Value = 3496.2563
$ws.Range("I132").Value = 2863.923
$ws.Range("K132").Value = 8591.769
$ws.Range("M132").Value = -6061.769
$ws.Range("H136").Value = 2585.0715
$ws.Range("I136").Value = 1954.7273
$ws.Range("J136").Value = 4896.3335
$ws.Range("K136").Value = 5864.1819
$ws.Range("L136").Value = 14689.0005
$ws.Range("M136").Value = -3314.1819
$ws.Range("N136").Value = -19789.0005

